# Weekly update: a new price record for "Acelga" (Feria Lagunitas de Puerto
# Montt) is inserted as the new row 31, pushing all subsequent rows (old
# 31..124) down by one (new 32..125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 31 - this shifts rows 31:124 down to 32:125
# and grows the used range to A1:R125 automatically.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record's data.
$ws.Cells.Item(31, 1).Value = 4
$ws.Cells.Item(31, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(31, 3).Value = "Los Lagos"
$ws.Cells.Item(31, 4).Value = 44525
$ws.Cells.Item(31, 5).Value = 10
$ws.Cells.Item(31, 6).Value = 100112009
$ws.Cells.Item(31, 7).Value = "Acelga"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 100
$ws.Cells.Item(31, 11).Value = 3500
$ws.Cells.Item(31, 12).Value = 3500
$ws.Cells.Item(31, 13).Value = 3500
$ws.Cells.Item(31, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(31, 15).Value = "Región del Maule"
$ws.Cells.Item(31, 16).Value = 875
$ws.Cells.Item(31, 17).Value = 4
$ws.Cells.Item(31, 18).Value = "Hortaliza"
